$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column ("TabName") shifting query/StatQuery/dbExcel/WebExcel
# columns one slot to the right.
$ws.Columns.Item(1).Insert()

# New header + value for the inserted column.
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Updated Cypher query text for the "query" column (now column B).
$casesQuery = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "ASIAN"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@
$ws.Range("B2").Value = $casesQuery

# Updated Cypher query text for the "StatQuery" column (now column C).
$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "ASIAN"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@
$ws.Range("C2").Value = $statQuery

# New column A is narrow (best-fit style width); keep the other columns as-is.
$ws.Columns.Item(1).ColumnWidth = 8

# Row 2 grows to fit the much longer, multi-line query text.
$ws.Rows.Item(2).RowHeight = 174

# Move the active selection like the author's saved view (now on the query cell).
$ws.Range("B2").Select() | Out-Null
